# Auto-generated edit script: fills in CrossRef-sourced metadata (Title, Abstract,
# Authors, ID Format, Date Accepted) for rows that previously had placeholder
# "not found" values, and refreshes the Authors column for rows that already had
# metadata (re-padded whitespace in the % separated author list), matching the
# "added and updated with crossref" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = @'
"Heterogeneous social interactions and the COVID-19 lockdown outcome in a multi-group SEIR model"
'@

$ws.Range("D2").Value = @'
"We study variants of the SEIR model for interpreting some qualitative features of the statistics of the Covid-19 epidemic in France.
 Standard SEIR models distinguish essentially two regimes: either the disease is controlled and the number of infected people rapidly decreases, or the disease spreads and contaminates a significant fraction of the population until herd immunity is achieved.
 After lockdown, at first sight it seems that social distancing is not enough to control the outbreak.
 We discuss here a possible explanation, namely that the lockdown is creating social heterogeneity: even if a large majority of the population complies with the lockdown rules, a small fraction of the population still has to maintain a normal or high level of social interactions, such as health workers, providers of essential services, etc.
 This results in an apparent high level of epidemic propagation as measured through re-estimations of the basic reproduction ratio.
 However, these measures are limited to averages, while variance inside the population plays an essential role on the peak and the size of the epidemic outbreak and tends to lower these two indicators.
 We provide theoretical and numerical results to sustain such a view.

'@

$ws.Range("E2").Value = @'
[Jean%Dolbeault%xref no email%1,  Gabriel%Turinici%xref no email%1]
'@

$ws.Range("F2").Value = @'
10.1051/mmnp/2020025
'@

$ws.Range("G2").Value = @'
CROSSREF
'@

$ws.Range("H2").Value = @'
2023-05-02
'@

$ws.Range("C3").Value = @'
"Early detection of superspreaders by mass group pool testing can mitigate COVID-19 pandemic"
'@

$ws.Range("D3").Value = @'
"SummaryBackgroundMost of epidemiological models applied for COVID-19 do not consider heterogeneity in infectiousness and impact of superspreaders, despite the broad viral loading distributions amongst COVID-19 positive people (1 – 106 per mL).
 Also, mass group testing is not used regardless to existing shortage of tests.
 I propose new strategy for early detection of superspreaders with reasonable number of RT-PCR tests, which can dramatically mitigate development COVID-19 pandemic and even turn it endemic.
MethodsI used stochastic social-epidemiological SEIAR model, where S-suspected, E-exposed, I-infectious, A-admitted (confirmed COVID-19 positive, who are admitted to hospital or completely isolated), R-recovered.
 The model was applied to real COVID-19 dynamics in London, Moscow and New York City.
FindingsViral loading data measured by RT-PCR were fitted by broad log-normal distribution, which governed high importance of superspreaders.
 The proposed full scale model of a metropolis shows that top 10% spreaders (100+ higher viral loading than median infector) transmit 45% of new cases.
 Rapid isolation of superspreaders leads to 4-8 fold mitigation of pandemic depending on applied quarantine strength and amount of currently infected people.
 High viral loading allows efficient group “matrix” pool testing of population focused on detection of the superspreaders requiring remarkably small amount of tests.
InterpretationThe model and new testing strategy may prevent thousand or millions COVID-19 deaths requiring just about 5000 daily RT-PCR test for big 12 million city such as Moscow.
 Though applied to COVID-19 pandemic the results are universal and can be used for other infectious heterogenous epidemics.
FundingNo funding
'@

$ws.Range("E3").Value = @'
[M.B.%Gongalsky%xref no email%1]
'@

$ws.Range("F3").Value = @'
10.1101/2020.04.22.20076166
'@

$ws.Range("G3").Value = @'
CROSSREF
'@

$ws.Range("H3").Value = @'
2022-11-21
'@

$ws.Range("E4").Value = @'
[ Nicolas%Hoertel%null%1,        Martin%Blachier%null%1,        Carlos%Blanco%null%1,        Mark%Olfson%null%1,        Marc%Massetti%null%1,        Marina Sánchez%Rico%null%1,        Frédéric%Limosin%null%1,        Henri%Leleu%null%1,       Nicolas%Hoertel%null%1,       Martin%Blachier%null%1,       Carlos%Blanco%null%1,       Mark%Olfson%null%1,       Marc%Massetti%null%1,       Marina Sánchez%Rico%null%1,       Frédéric%Limosin%null%1,       Henri%Leleu%null%1]
'@

$ws.Range("E5").Value = @'
[Tuo%Ji%NULL%1,       Hai-Lian%Chen%NULL%1,       Jing%Xu%NULL%1,       Ling-Ning%Wu%NULL%1,       Jie-Jia%Li%NULL%1,       Kai%Chen%NULL%1,       Gang%Qin%tonygqin@ntu.edu.cn%1]
'@

$ws.Range("C6").Value = @'
"The effectiveness of quarantine of Wuhan city against the Corona Virus Disease 2019 (COVID‐19): A well‐mixed SEIR model analysis"
'@

$ws.Range("E6").Value = @'
[Can%Hou%xref no email%1,  Jiaxin%Chen%xref no email%1,  Yaqing%Zhou%xref no email%1,  Lei%Hua%xref no email%1,  Jinxia%Yuan%xref no email%1,  Shu%He%xref no email%1,  Yi%Guo%xref no email%1,  Sheng%Zhang%xref no email%2,  Qiaowei%Jia%xref no email%1,  Chenhui%Zhao%xref no email%1,  Jing%Zhang%xref no email%1,  Guangxu%Xu%xref no email%1,  Enzhi%Jia%xref no email%1]
'@

$ws.Range("F6").Value = @'
10.1002/jmv.25827
'@

$ws.Range("G6").Value = @'
CROSSREF
'@

$ws.Range("H6").Value = @'
2023-05-23
'@

$ws.Range("E7").Value = @'
[Hien%Lau%hlau2@uci.edu%1,       Veria%Khosrawipour%veriakhosrawipour@yahoo.de%1,       Piotr%Kocbach%piotr.kocbach@uwm.edu.pl%1,       Agata%Mikolajczyk%agata.mikolajczyk@upwr.edu.pl%1,       Justyna%Schubert%justyna.schubert@upwr.edu.pl%1,       Jacek%Bania%jacek.bania@upwr.edu.pl%1,       Tanja%Khosrawipour%tkhosrawipour@gmail.com%1]
'@

$ws.Range("E8").Value = @'
[Carlo%Signorelli%NULL%1,       Thea%Scognamiglio%NULL%1,       Anna%Odone%NULL%1]
'@

$ws.Range("C9").Value = @'
"Population-scale testing can suppress the spread of COVID-19"
'@

$ws.Range("D9").Value = @'
"We propose an additional intervention that would contribute to the control of the COVID-19 pandemic, offer more protection for people working in essential jobs, and help guide an eventual reopening of society.
 The intervention is based on: (1) testing every individual (2) repeatedly, and (3) self-quarantine of infected individuals.
 Using a standard epidemiological model (SIR), we show here that by identification and isolation of the majority of infectious individuals, including those who may be asymptomatic, the reproduction number R0of SARS-CoV-2 would be reduced well below 1.0, and the epidemic would collapse.
 We replicate these observations in a more complex stochastic dynamic model on a social network graph.
 We also find that the testing regime would be additive to other interventions, and be effective at any level of prevalence.
 If adopted as a policy, any industrial society could sustain the regime for as long as it takes to find a safe and effective cure or vaccine.
 Our model also indicates that unlike sampling-based tests, population-scale testing does not need to be very accurate: false negative rates up to 15% could be tolerated if 80% comply with testing every ten days, and false positives can be almost arbitrarily high when a high fraction of the population is already effectively quarantined.
 Testing at the required scale would be feasible if existing qPCR-based methods are scaled up and multiplexed.
 A mass produced, low throughput field test kit could also be carried out at home.
 Economic analysis also supports the feasibility of the approach: current reagent costs for tests are in the range of a dollar or less, and the estimated benefits for population-scale testing are so large that the policy would be cost-effective even if the costs were larger by more than two orders of magnitude.
 To identify both active and previous infections, both viral RNA and antibodies could be tested.
 All technologies to build such test kits, and to produce them in the scale required to test the entire worlds’ population exist already.
 Integrating them, scaling up production, and implementing the testing regime will require resources and planning, but at a scale that is very small compared to the effort that every nation would devote to defending itself against a more traditional foe.

'@

$ws.Range("E9").Value = @'
[Jussi%Taipale%xref no email%1,  Paul%Romer%xref no email%1,  Sten%Linnarsson%xref no email%1]
'@

$ws.Range("F9").Value = @'
10.1101/2020.04.27.20078329
'@

$ws.Range("G9").Value = @'
CROSSREF
'@

$ws.Range("H9").Value = @'
2023-04-28
'@

$ws.Range("E11").Value = @'
[Nadya%Johanna%NULL%1,       Henrico%Citrawijaya%NULL%1,       Grace%Wangge%NULL%1]
'@

$ws.Range("C12").Value = @'
"Price of Delay in Covid-19 Lockdowns: Delays Spike Total Cases, Natural Experiments Reveal"
'@

$ws.Range("E12").Value = @'
[Gerard J.%Tellis%xref no email%1,  Ashish%Sood%xref no email%1,  Nitish%Sood%xref no email%1]
'@

$ws.Range("F12").Value = @'
10.2139/ssrn.3592912
'@

$ws.Range("G12").Value = @'
CROSSREF
'@

$ws.Range("H12").Value = @'
2022-04-27
'@

$ws.Range("E13").Value = @'
[Biao%Tang%NULL%1,       Fan%Xia%NULL%1,       Sanyi%Tang%NULL%1,       Nicola Luigi%Bragazzi%NULL%1,       Qian%Li%NULL%1,       Xiaodan%Sun%NULL%1,       Juhua%Liang%NULL%1,       Yanni%Xiao%yxiao@mail.xjtu.edu.cn%1,       Jianhong%Wu%wujh@yorku.ca%1]
'@
